# Apply the "Manual reviewed" update to the effort sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effort")

# New row 33: 2012-10-30, 2.75h, 2h additional, "Manual continued"
$ws.Cells.Item(33, 1).Value = 41212
$ws.Cells.Item(33, 2).Value = 2.75
$ws.Cells.Item(33, 3).Value = 2
$ws.Cells.Item(33, 4).Value = "Manual continued"

# New row 34: 2012-10-31, 2.75h, "Manual reviewed"
$ws.Cells.Item(34, 1).Value = 41213
$ws.Cells.Item(34, 2).Value = 2.75
$ws.Cells.Item(34, 4).Value = "Manual reviewed"

# New row 35: 2012-11-01, 1.5h, "Manual reviewed"
$ws.Cells.Item(35, 1).Value = 41214
$ws.Cells.Item(35, 2).Value = 1.5
$ws.Cells.Item(35, 4).Value = "Manual reviewed"

# Update the view: scroll and select like in the target
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("D35").Select()
